# Rename the embedded logo pictures referenced from the document's
# headers/footers, per the authored commit:
#   - the two Pearson Edexcel logo pictures (in the "first page" footer
#     and the "default" footer) are renamed from image1.png -> image2.png
#   - the BTec_Logo-Orange picture (in the "first page" header) is
#     renamed from image2.jpg -> image1.jpg
#
# These "name" values live on the inline picture's docPr/cNvPr pair; the
# Word object model exposes them through InlineShape.Name (after
# converting to a floating Shape, which is the supported way to read the
# current value back), so walk every header/footer of every section and
# rename whichever picture we find by matching its stable identity
# (AlternativeText, i.e. the picture's "descr").

$d = $word.ActiveDocument

function Rename-InlineShapeByAltText($range, $altText, $newName) {
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.AlternativeText -eq $altText) {
            $shp.Name = $newName
        }
    }
}

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)

    for ($hi = 1; $hi -le 3; $hi++) {
        $hf = $sec.Headers.Item($hi)
        if ($hf.Exists) {
            Rename-InlineShapeByAltText $hf.Range "BTec_Logo-Orange" "image1.jpg"
        }
    }

    for ($fi = 1; $fi -le 3; $fi++) {
        $ft = $sec.Footers.Item($fi)
        if ($ft.Exists) {
            Rename-InlineShapeByAltText $ft.Range "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" "image2.png"
        }
    }
}
